# Adds 2022-Q1 data.
#
# The workbook's last sheet was "总计" (grand-totals-by-quarter). We:
#   1. Duplicate that sheet (Worksheet.Copy) so both copies start out with
#      identical sheet-level properties (sheetPr/pageMargins/etc.) to the
#      original "总计" sheet.
#   2. Rename the original to "2022-Q1" (keeps its original sheetId/r:id)
#      and repurpose it to hold the new quarter's per-fund holding detail
#      table (same shape as the other "20XX-QX" sheets).
#   3. Rename the duplicate to "总计" (it gets a fresh sheetId/r:id) and
#      repopulate it with the refreshed totals-by-quarter table: the new
#      2022-Q1 summary row on top, followed by all the previously
#      existing quarters.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# A sheet that already has the header/index cell formatting (style s="2":
# bold font + thin border + centered alignment) we want to reuse without
# minting any new style entries.
$styleSource = $wb.Worksheets.Item("2021-Q4")
$headerStyleCell = $styleSource.Range("B1")
$indexStyleCell  = $styleSource.Range("A2")

$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Copy($null, $oldTotal)

$q1 = $wb.Worksheets.Item("总计")
$total = $wb.Worksheets.Item("总计 (2)")

# ---------------------------------------------------------------------
# Step 1: Turn the original "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerStyleCell.Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Data rows (fund code, fund name, fund size, total stock position,
# position ratio, held market value, position rank). Per the source data,
# columns B and D-G are text (fund codes keep leading zeros, and the
# numeric-looking values keep a fixed number of decimal places), while
# columns A and H are real numbers.
$fundRows = @(
    @("011184", "东方阿尔法招阳混合A", "6.40", "92.60", "9.50", "0.6080", 1),
    @("012368", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金A", "5.76", "93.78", "9.41", "0.5420", 2),
    @("005358", "东方阿尔法精选灵活配置混合A", "4.31", "92.94", "9.27", "0.3995", 1),
    @("233006", "大摩领先优势混合", "4.12", "94.42", "9.30", "0.3832", 2),
    @("000309", "大摩品质生活精选股票", "4.36", "94.17", "8.31", "0.3623", 3),
    @("010322", "摩根士丹利华鑫新兴产业股票", "2.41", "94.11", "9.63", "0.2321", 2),
    @("002707", "摩根士丹利华鑫科技领先灵活配置混合", "2.27", "93.05", "7.13", "0.1619", 2),
    @("005359", "东方阿尔法精选灵活配置混合C", "0.54", "92.94", "9.27", "0.0501", 1),
    @("012369", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金C", "0.40", "93.78", "9.41", "0.0376", 2),
    @("011185", "东方阿尔法招阳混合C", "0.08", "92.60", "9.50", "0.0076", 1)
)

# Force columns B and D-G to Text *before* writing, so fund codes such as
# "011184" and decimal strings such as "6.40" / "0.6080" are stored
# verbatim instead of being auto-converted to numbers.
$q1.Range("B2:B11").NumberFormat = "@"
$q1.Range("D2:G11").NumberFormat = "@"

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A$r").Value = ($r - 2)
    $q1.Range("B$r").Value = $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = $row[2]
    $q1.Range("E$r").Value = $row[3]
    $q1.Range("F$r").Value = $row[4]
    $q1.Range("G$r").Value = $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# Drop the forced-text direct formatting again (back to the default/
# unstyled cell) without touching the values/types that were already
# written, then stamp column A with the shared index-cell style.
$q1.Range("B2:G11").Style = "Normal"
$indexStyleCell.Copy()
$q1.Range("A2:A11").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Step 2: Turn the duplicated sheet into the refreshed "总计" sheet
# ---------------------------------------------------------------------
$total.Name = "总计"
$total.Cells.Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$headerStyleCell.Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)

$totalRows = @(
    @("2022-Q1", 10, 2.78),
    @("2021-Q4", 10, 3.57),
    @("2021-Q3", 10, 3.12),
    @("2021-Q2", 11, 3.79),
    @("2021-Q1", 14, 2.85),
    @("2020-Q4", 9, 2.54)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = ($r - 2)
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$indexStyleCell.Copy()
$total.Range("A2:A7").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Leave the originally active sheet ("2020-Q4") selected/active, matching
# the workbook's prior state.
$wb.Worksheets.Item("2020-Q4").Activate()
